$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (D) and Volume(1h) (E) values from the refreshed
# coinranking.com symbol-list snapshot. Values are written as literal text
# (leading apostrophe) so Excel does not reinterpret the numeric-looking
# strings/percentages as numbers, then ClearFormats() drops the resulting
# quote-prefix formatting so the cell style stays the original "General" one.
$updates = @(
    @("D2", "308.50"),
    @("E2", "3.70%"),
    @("D3", "44.28"),
    @("E3", "6.77%"),
    @("D4", "5.087"),
    @("E4", "1.18%"),
    @("D5", "0.07997"),
    @("E5", "5.91%"),
    @("D6", "4.454"),
    @("E6", "1.82%"),
    @("D7", "1.616"),
    @("E7", "1.30%"),
    @("D8", "1.068"),
    @("E8", "15.09%"),
    @("D9", "0.1289"),
    @("E9", "7.65%"),
    @("E10", "2.56%"),
    @("D11", "0.09236"),
    @("E11", "3.93%"),
    @("D12", "0.04179"),
    @("E12", "3.50%"),
    @("D13", "0.1035"),
    @("E13", "-1.76%"),
    @("D14", "0.001297"),
    @("D15", "0.005694"),
    @("E15", "-1.85%"),
    @("D17", "3.366"),
    @("E17", "0.81%"),
    @("D18", "2.399"),
    @("E18", "-0.19%"),
    @("D19", "0.3352"),
    @("E19", "1.22%"),
    @("D20", "8.000"),
    @("E20", "0.27%"),
    @("D21", "0.1372"),
    @("E21", "-3.25%"),
    @("D22", "0.3115"),
    @("E22", "3.89%"),
    @("D23", "0.04146"),
    @("E23", "2.36%"),
    @("D24", "0.001262"),
    @("E24", "-0.31%"),
    @("D25", "0.004310"),
    @("E25", "3.48%"),
    @("D26", "0.0001331"),
    @("E26", "8.28%"),
    @("D38", "0.02654"),
    @("E38", "9.64%"),
    @("D39", "0.05376"),
    @("E39", "3.18%"),
    @("D40", "0.005589"),
    @("E40", "-14.64%"),
    @("D41", "0.007714"),
    @("E41", "-1.04%"),
    @("E42", "5.51%"),
    @("D43", "0.007247"),
    @("E43", "-4.07%"),
    @("D44", "0.008394"),
    @("E44", "7.20%"),
    @("D45", "0.3079"),
    @("E45", "-4.51%"),
    @("D46", "0.00006686"),
    @("E46", "-1.46%"),
    @("D47", "0.00000000740"),
    @("E47", "-1.34%"),
    @("D48", "0.06162"),
    @("E48", "33.18%"),
    @("D49", "0.003944"),
    @("E49", "-6.04%"),
    @("D50", "0.00002071"),
    @("E50", "-1.34%"),
    @("D51", "0.0001972"),
    @("E51", "-1.34%"),
)

foreach ($update in $updates) {
    $cellRef = $update[0]
    $newValue = $update[1]
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $newValue
    $range.ClearFormats()
}
